# Apply the "ginger-comments-2013-07-23" review-response edit:
#   - Row 12 (project 20120232): reviewer adds "OK" in column L, replying to
#     the "multiple points per project ok? - YES" comment in K12.
#   - Row 13 (project 20120272): reviewer adds "OK, fixed." in column L,
#     replying to the misidentified-project comment in K13.
#   - Row 14 (project 20120284): reviewer adds "OK, done." in column L
#     (same stock response already used elsewhere in the sheet, e.g. L7/L11).
#   - Rows 13/14 shrink slightly in height to fit the now-shorter wrapped text.
#   - The active selection moves from L11 to L14, with the view scrolled so
#     row 12 is at/near the top.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L12").Value = "OK"
$ws.Range("L13").Value = "OK, fixed."
$ws.Range("L14").Value = "OK, done."

$ws.Rows.Item(13).RowHeight = 55.2
$ws.Rows.Item(14).RowHeight = 149.25

# Move the viewport/selection to match the post-edit state (top-left A12,
# active cell L14).
$aw = $excel.ActiveWindow
$aw.ScrollRow = 12
$aw.ScrollColumn = 1
$ws.Range("L14").Select()
